$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.809.10'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.624.85'
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.05'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  -0.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.60'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0786'
$ws.Range('E11').Value = '  -0.78%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.851.67'
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.22'
$ws.Range('E13').Value = '  -0.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.623.40'
$ws.Range('E14').Value = '  -0.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.541'
$ws.Range('E15').Value = '  -2.53%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₃0755'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.55'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.805.02'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('E19').Value = '  -0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.97'
$ws.Range('E20').Value = '  +0.30%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.35'
$ws.Range('E21').Value = '  -1.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.93'
$ws.Range('E22').Value = '  +0.38%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.22'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('E24').Value = '  -2.35%  '
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.32'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('E27').Value = '  +0.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.84'
$ws.Range('E28').Value = '  -0.28%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.42'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.23'
$ws.Range('E30').Value = '  -0.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0497'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.30'
$ws.Range('E32').Value = '  -0.88%  '
$ws.Range('E33').Value = '  -0.36%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  +2.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.899'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.127.97'
$ws.Range('E37').Value = '  -0.51%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.545'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('E39').Value = '  -2.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0155'
$ws.Range('E40').Value = '  +0.53%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.45'
$ws.Range('E42').Value = '  -1.08%  '
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.06'
$ws.Range('E43').Value = '  -1.31%  '
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.762.24'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0₆0111'
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.97'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0530'
$ws.Range('E48').Value = '  +4.53%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.44'
$ws.Range('E49').Value = '  +1.87%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.414'
$ws.Range('E50').Value = '  -1.02%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.60'
$ws.Range('E51').Value = '  +2.70%  '
